$d = $word.ActiveDocument
$paras = $d.Paragraphs
$target = $null
for ($i = 1; $i -le $paras.Count; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    if (($t -like "*Certified that the above*") -and ($t -like "*PassportNumber*")) {
        $target = $p.Range
    }
}
if ($target -eq $null) {
    Write-Host "ERROR: target paragraph not found"
} else {
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="004F1B8E" w:rsidRPr="008E6C50" w:rsidRDefault="00346247" w:rsidP="006C4C85"><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="both"/><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr><w:r w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">              </w:t></w:r><w:r w:rsidR="005C4E0A" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">  </w:t></w:r><w:proofErr w:type="gramStart"/><w:r w:rsidR="004F1B8E" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Certified that the above</w:t></w:r><w:r w:rsidR="00D555ED" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> impressions taken by me are of</w:t></w:r><w:r w:rsidR="007E2ADB" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="007E2ADB" w:rsidRPr="008E6C50"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>M</w:t></w:r><w:r w:rsidR="00B371D9" w:rsidRPr="008E6C50"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>r</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>/M</w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">s. </w:t></w:r><w:bookmarkStart w:id="6" w:name="NameAddress"/><w:r w:rsidR="00B371D9" w:rsidRPr="008E6C50"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Name</w:t></w:r><w:bookmarkEnd w:id="6"/><w:r w:rsidR="00D555ED" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r w:rsidR="004F1B8E" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">holder of </w:t></w:r><w:r w:rsidR="00176EC4" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Indian Passport</w:t></w:r><w:r w:rsidR="004F1B8E" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> No.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:bookmarkStart w:id="7" w:name="PassportNo"/><w:proofErr w:type="spellStart"/><w:proofErr w:type="gramStart"/><w:r w:rsidR="00B371D9" w:rsidRPr="008E6C50"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Passport</w:t></w:r><w:r w:rsidR="00B371D9" w:rsidRPr="008E6C50"><w:rPr><w:b/><w:bCs/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Number</w:t></w:r><w:bookmarkEnd w:id="7"/><w:proofErr w:type="spellEnd"/><w:r w:rsidR="007E2ADB" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r w:rsidR="004F1B8E" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">  The impressions are taken for the purpose of </w:t></w:r><w:r w:rsidR="00DC3B48" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">applying for </w:t></w:r><w:bookmarkStart w:id="8" w:name="Reason"/><w:r w:rsidR="00B371D9" w:rsidRPr="008E6C50"><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>Reason</w:t></w:r><w:bookmarkEnd w:id="8"/><w:r w:rsidR="007E2ADB" w:rsidRPr="008E6C50"><w:rPr><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>.</w:t></w:r></w:p>'
    $ret = $target.InsertXML($xml)
    Write-Host "InsertXML ret=$ret"
}
